# Update countries & provincias Spain
# Applies the scraped-data refresh: reorders a couple of country rows
# (their per-day numbers shift rank), updates several countries' case
# counters, and bumps the "last updated" timestamp string.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Timestamp string (A1) -------------------------------------------------
$ws.Range("A1").Value = "Datos actualizados a 20 de Julio de 2020 a las 04:31"

# --- Country ranking shuffle: rows 37-39 (Kuwait/Ucrania/Bolivia) ----------
# New data lands in row 37 (now Bolivia); the old row37/row38 data rolls
# down into row38/row39 (now Kuwait / Ucrania respectively).
$ws.Range("A37").Value = "Bolivia"
$ws.Range("B37").Value = 59582
$ws.Range("C37").Value = 1444
$ws.Range("D37").Value = 18553
$ws.Range("E37").Value = 38878
$ws.Range("F37").Value = 0
$ws.Range("G37").Value = 45
$ws.Range("H37").Value = 2151

$ws.Range("A38").Value = "Kuwait"
$ws.Range("B38").Value = 59204
$ws.Range("C38").Value = 0
$ws.Range("D38").Value = 49687
$ws.Range("E38").Value = 9109
$ws.Range("F38").Value = 0
$ws.Range("G38").Value = 0
$ws.Range("H38").Value = 408

$ws.Range("A39").Value = "Ucrania"
$ws.Range("B39").Value = 58842
$ws.Range("C39").Value = 0
$ws.Range("D39").Value = 30879
$ws.Range("E39").Value = 26478
$ws.Range("F39").Value = 0
$ws.Range("G39").Value = 0
$ws.Range("H39").Value = 1485

# --- Country ranking swap: rows 210-211 (Islas Malvinas / Groenlandia) -----
# Identical totals (tied), so only the labels trade places.
$ws.Range("A210").Value = "Groenlandia"
$ws.Range("A211").Value = "Islas Malvinas"

# --- Mexico (row 10) updated counters --------------------------------------
$ws.Range("D10").Value = 217423
$ws.Range("E10").Value = 87617
$ws.Range("G10").Value = 296
$ws.Range("H10").Value = 39184

# --- China (row 29) updated counters ----------------------------------------
$ws.Range("B29").Value = 83682
$ws.Range("C29").Value = 22
$ws.Range("D29").Value = 78799
$ws.Range("E29").Value = 249

# --- Corea del Sur (row 71) updated counters --------------------------------
$ws.Range("B71").Value = 13771
$ws.Range("C71").Value = 26
$ws.Range("D71").Value = 12572
$ws.Range("E71").Value = 903
$ws.Range("G71").Value = 1
$ws.Range("H71").Value = 296

# --- Nueva Zelanda (row 134) updated counters -------------------------------
$ws.Range("B134").Value = 1554
$ws.Range("C134").Value = 1
$ws.Range("E134").Value = 26
